$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in self-assessment scores that were previously left blank.
# I24 already shares formatting with its column neighbours, so a plain
# value write is enough.
$ws.Range("I24").Value = 5

# E27:E29 were blank and therefore still carrying the "no value yet"
# left/wrap style; once a score is entered they should pick up the same
# centered numeric style used by the rest of the "自评分" column (e.g. E8).
# Cloning formats from a sibling score cell reproduces that style switch.
$ws.Range("E8").Copy()
$ws.Range("E27:E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E27").Value = 4
$ws.Range("E28").Value = 4
$ws.Range("E29").Value = 4

# Move the active selection down to the bottom of the form (rows 33-38),
# matching the reviewer's current scroll position.
$ws.Range("E32").Select()
